$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Ligand/Receptor derived specificity values updated with new TPM calc)
$ws.Range("I2").Value = 0.9891011365778073
$ws.Range("J2").Value = 0.9927075980877177
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("S2").Value = 0.9891011365778073
$ws.Range("T2").Value = 0.9927075980877177

# Row 3 updates: Sending/Target clusters swapped (FAPs -> MuSCs, MuSCs -> FAPs)
$ws.Range("A3").Value = "MuSCs"
$ws.Range("D3").Value = "FAPs"

# Row 3 numeric values recomputed with new TPM data
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.0206735
$ws.Range("H3").Value = 0.041347
$ws.Range("I3").Value = 0.01089886342219268
$ws.Range("J3").Value = 0.007292401912282354
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1753453333333333
$ws.Range("N3").Value = 0.5260359999999999
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.003625001748666667
$ws.Range("R3").Value = 0.021750010492
$ws.Range("S3").Value = 0.01089886342219268
$ws.Range("T3").Value = 0.007292401912282354

$wb.Save()
